$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each (cell, new value) pair below reproduces the 15-Jan-2023 data refresh
# described in the commit. Values are written as literal text (NumberFormat
# forced to "@" beforehand) so price/percentage strings such as "-2.49%" are
# preserved verbatim instead of being reinterpreted as numbers/percentages by
# Excel's normal type inference, then the style is reset to "Normal" so no
# stray number-format is left behind on the cell.

$updates = @(
    @{ Cell = 'D2'; Value = '298.99' }
    @{ Cell = 'E2'; Value = '-2.49%' }
    @{ Cell = 'D3'; Value = '31.75' }
    @{ Cell = 'E3'; Value = '-1.44%' }
    @{ Cell = 'D4'; Value = '5.124' }
    @{ Cell = 'E4'; Value = '-3.82%' }
    @{ Cell = 'D5'; Value = '0.07526' }
    @{ Cell = 'E5'; Value = '0.87%' }
    @{ Cell = 'D6'; Value = '7.787' }
    @{ Cell = 'E6'; Value = '0.48%' }
    @{ Cell = 'D7'; Value = '1.716' }
    @{ Cell = 'E7'; Value = '10.26%' }
    @{ Cell = 'D8'; Value = '3.796' }
    @{ Cell = 'E8'; Value = '2.52%' }
    @{ Cell = 'D9'; Value = '0.9243' }
    @{ Cell = 'E9'; Value = '0.17%' }
    @{ Cell = 'D10'; Value = '0.1711' }
    @{ Cell = 'E10'; Value = '2.75%' }
    @{ Cell = 'D11'; Value = '0.07387' }
    @{ Cell = 'E11'; Value = '-2.54%' }
    @{ Cell = 'D12'; Value = '0.07946' }
    @{ Cell = 'E12'; Value = '-0.16%' }
    @{ Cell = 'D13'; Value = '0.03040' }
    @{ Cell = 'E13'; Value = '-1.09%' }
    @{ Cell = 'D14'; Value = '0.09901' }
    @{ Cell = 'E14'; Value = '0.58%' }
    @{ Cell = 'D15'; Value = '0.001503' }
    @{ Cell = 'E15'; Value = '-1.38%' }
    @{ Cell = 'D16'; Value = '0.04658' }
    @{ Cell = 'E16'; Value = '2.41%' }
    @{ Cell = 'D17'; Value = '0.006233' }
    @{ Cell = 'E17'; Value = '-2.55%' }
    @{ Cell = 'D18'; Value = '3.451' }
    @{ Cell = 'E18'; Value = '-0.74%' }
    @{ Cell = 'D19'; Value = '2.220' }
    @{ Cell = 'E19'; Value = '-0.82%' }
    @{ Cell = 'D20'; Value = '0.3292' }
    @{ Cell = 'E20'; Value = '0.44%' }
    @{ Cell = 'E21'; Value = '0.72%' }
    @{ Cell = 'D22'; Value = '4.559' }
    @{ Cell = 'E22'; Value = '8.32%' }
    @{ Cell = 'D23'; Value = '0.1550' }
    @{ Cell = 'E23'; Value = '-4.74%' }
    @{ Cell = 'D24'; Value = '0.001218' }
    @{ Cell = 'D25'; Value = '0.004422' }
    @{ Cell = 'E25'; Value = '-2.39%' }
    @{ Cell = 'D26'; Value = '0.0001399' }
    @{ Cell = 'E26'; Value = '19.79%' }
    @{ Cell = 'D27'; Value = '0.0001843' }
    @{ Cell = 'E27'; Value = '10.72%' }
    @{ Cell = 'D39'; Value = '0.01666' }
    @{ Cell = 'E39'; Value = '0.31%' }
    @{ Cell = 'D40'; Value = '0.04554' }
    @{ Cell = 'E40'; Value = '0.79%' }
    @{ Cell = 'D41'; Value = '0.007115' }
    @{ Cell = 'E41'; Value = '-4.01%' }
    @{ Cell = 'E42'; Value = '-2.80%' }
    @{ Cell = 'E43'; Value = '-8.75%' }
    @{ Cell = 'D44'; Value = '0.01285' }
    @{ Cell = 'E44'; Value = '-6.99%' }
    @{ Cell = 'D45'; Value = '0.00006062' }
    @{ Cell = 'E45'; Value = '0.64%' }
    @{ Cell = 'D46'; Value = '1.930' }
    @{ Cell = 'E46'; Value = '1.97%' }
    @{ Cell = 'E47'; Value = '-5.61%' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}

